# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '33.624.72'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.769.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.54%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.24'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.554'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.55%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.73'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.74'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.276'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0659'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0921'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.028.37'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.774.55'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.621'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.650.31'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '9.94'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.16'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '68.19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '248.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0734'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.20'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.15'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.14'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.69'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.90'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.62%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.78'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0513'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.53'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.81'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.473.84'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.06'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.624'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0183'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.51'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.36'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.879'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.06'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.45%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.08'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.88%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0507'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.922.21'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.95%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.68'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.60'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +12.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.50'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.07%  '
